# Se agrega validación de fecha
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base Clientes carga manual")
$lo = $ws.ListObjects.Item("Tabla1")

# New client rows to append to the "Tabla1" table on the sheet.
$newRows = @(
    @(10240057, "MARIA GUADALUPE RUIZ ESTRADA"),
    @(10174274, "AGRICOLA EL MORON SA DE CV"),
    @(500231,   "INSECTICIDAS HERBICIDAS Y SEMILLAS RIDA"),
    @(10181721, "DUNE COMPANY MEXICALI")
)

foreach ($rowData in $newRows) {
    $newRow = $lo.ListRows.Add()
    $newRow.Range.Cells.Item(1, 1).Value = $rowData[0]
    $newRow.Range.Cells.Item(1, 2).Value = $rowData[1]
}

# Keep column B wide enough to fit the longest client name.
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

$ws.Range("B14").Select() | Out-Null
